$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus results for the 380 kV case), columns B-F and I-N.
$data = @{
    2 = @{ "B" = 1.02; "C" = 1.078166994231856; "D" = 1.078379102296422; "E" = 1.081245715788358; "F" = 1.09137242642333; "I" = 1.052585734503189; "J" = 1.083057820201336; "K" = 1.081057929975439; "L" = 1.0839170409593; "M" = 1.094017588333633; "N" = 1.08459588653449 }
    3 = @{ "B" = 1.02; "C" = 1.079968708364879; "D" = 1.079823307657738; "E" = 1.082852628289828; "F" = 1.093022595092452; "I" = 1.053093621682712; "J" = 1.084516075177242; "K" = 1.082318480791172; "L" = 1.085340453596157; "M" = 1.09548607232205; "N" = 1.086056212399726 }
    4 = @{ "B" = 1.02; "C" = 1.081131555418659; "D" = 1.080755026700154; "E" = 1.083889449513674; "F" = 1.094087687258353; "I" = 1.053419387390923; "J" = 1.085456350196565; "K" = 1.083130841339589; "L" = 1.086258048708481; "M" = 1.096433116379492; "N" = 1.086997822717428 }
    5 = @{ "B" = 1.02; "C" = 1.081619718299505; "D" = 1.081146068093253; "E" = 1.084324634778359; "F" = 1.094534824732614; "I" = 1.053555657119054; "J" = 1.085850861629324; "K" = 1.083471578832369; "L" = 1.086642993257186; "M" = 1.096830509210465; "N" = 1.087392894401672 }
    6 = @{ "B" = 1.02; "C" = 1.081701642524894; "D" = 1.081211687697012; "E" = 1.08439766398278; "F" = 1.094609864693713; "I" = 1.05357849755326; "J" = 1.085917056461739; "K" = 1.083528744805396; "L" = 1.086707579875607; "M" = 1.096897189927976; "N" = 1.087459183238342 }
    7 = @{ "B" = 1.02; "C" = 1.081138080995855; "D" = 1.080760254365098; "E" = 1.083895267189848; "F" = 1.09409366437717; "I" = 1.053421210907009; "J" = 1.085461624725342; "K" = 1.08313539733579; "L" = 1.086263195530858; "M" = 1.096438429269741; "N" = 1.087003104736641 }
    8 = @{ "B" = 1.02; "C" = 1.078776519162145; "D" = 1.078867760717562; "E" = 1.081789400406827; "F" = 1.091930671274854; "I" = 1.052757975017812; "J" = 1.083551338269377; "K" = 1.08148462909329; "L" = 1.084398811742121; "M" = 1.094514532091073; "N" = 1.085090105454795 }
    9 = @{ "B" = 1.02; "C" = 1.07459158278698; "D" = 1.075511130816403; "E" = 1.078055307205318; "F" = 1.088098079346159; "I" = 1.051567042900392; "J" = 1.080159203284838; "K" = 1.078549985474744; "L" = 1.081086533978311; "M" = 1.091099575182507; "N" = 1.081693153249495 }
    10 = @{ "B" = 1.02; "C" = 1.071784771060945; "D" = 1.073257958245904; "E" = 1.075549373640774; "F" = 1.085527950281398; "I" = 1.050757813178409; "J" = 1.077879493977775; "K" = 1.076575515011339; "L" = 1.078859373448275; "M" = 1.088805448792877; "N" = 1.07941020649373 }
    11 = @{ "B" = 1.02; "C" = 1.070565166613561; "D" = 1.07227848366277; "E" = 1.074460156868794; "F" = 1.084411288048384; "I" = 1.050403710192859; "J" = 1.076887832096593; "K" = 1.075716109365099; "L" = 1.077890302921544; "M" = 1.087807736950089; "N" = 1.078417136338931 }
    12 = @{ "B" = 1.02; "C" = 1.070111495133703; "D" = 1.071914072227584; "E" = 1.074054935530638; "F" = 1.083995925337464; "I" = 1.050271618008571; "J" = 1.07651878765098; "K" = 1.075396205937239; "L" = 1.077529625542336; "M" = 1.087436474547505; "N" = 1.078047567807876 }
    13 = @{ "B" = 1.02; "C" = 1.070208839212594; "D" = 1.071992266696888; "E" = 1.074141886018488; "F" = 1.084085048777314; "I" = 1.050299977793302; "J" = 1.07659798074135; "K" = 1.075464857389582; "L" = 1.077607024970241; "M" = 1.087516142067673; "N" = 1.078126873361517 }
    14 = @{ "B" = 1.02; "C" = 1.070527679485691; "D" = 1.072248373440722; "E" = 1.074426674255704; "F" = 1.084376966061621; "I" = 1.050392802925887; "J" = 1.076857341081301; "K" = 1.075689680039107; "L" = 1.077860504027304; "M" = 1.087777061998898; "N" = 1.078386602022899 }
    15 = @{ "B" = 1.02; "C" = 1.070724039972881; "D" = 1.072406090508446; "E" = 1.074602056667416; "F" = 1.084556747957927; "I" = 1.050449920811692; "J" = 1.077017048693951; "K" = 1.075828109918471; "L" = 1.078016584760337; "M" = 1.087937734414936; "N" = 1.078546536438678 }
    16 = @{ "B" = 1.02; "C" = 1.071865621265977; "D" = 1.073322880681826; "E" = 1.07562157288051; "F" = 1.085601978416934; "I" = 1.050781235294098; "J" = 1.077945210443613; "K" = 1.076632456067649; "L" = 1.078923587179462; "M" = 1.088871570861945; "N" = 1.079476016284487 }
    17 = @{ "B" = 1.02; "C" = 1.072580557062695; "D" = 1.073896921645364; "E" = 1.076259970471776; "F" = 1.08625660020128; "I" = 1.050988064636122; "J" = 1.078526196959114; "K" = 1.077135800896989; "L" = 1.079491257512373; "M" = 1.089456168784055; "N" = 1.08005782786748 }
    18 = @{ "B" = 1.02; "C" = 1.072997159597032; "D" = 1.074231380908437; "E" = 1.076631939568737; "F" = 1.086638066503293; "I" = 1.051108348026244; "J" = 1.078864640319299; "K" = 1.077428965299293; "L" = 1.079821918074281; "M" = 1.089796737191033; "N" = 1.080396751856057 }
    19 = @{ "B" = 1.02; "C" = 1.073139141779588; "D" = 1.074345360705476; "E" = 1.076758704523686; "F" = 1.086768075502946; "I" = 1.051149301252806; "J" = 1.07897996724797; "K" = 1.07752885458296; "L" = 1.079934588606048; "M" = 1.089912791846871; "N" = 1.080512242562194 }
    20 = @{ "B" = 1.02; "C" = 1.072503893500643; "D" = 1.073835370768076; "E" = 1.076191517682068; "F" = 1.086186403200036; "I" = 1.050965910743725; "J" = 1.078463907843293; "K" = 1.07708184114259; "L" = 1.079430398719811; "M" = 1.089393490274995; "N" = 1.07999545029397 }
    21 = @{ "B" = 1.02; "C" = 1.070433807241897; "D" = 1.072172972819975; "E" = 1.074342828960785; "F" = 1.084291019972235; "I" = 1.050365483820285; "J" = 1.076780985298174; "K" = 1.075623494308034; "L" = 1.077785880807569; "M" = 1.087700246144541; "N" = 1.078310137805803 }
    22 = @{ "B" = 1.02; "C" = 1.069128456065898; "D" = 1.071124330730154; "E" = 1.073176786658082; "F" = 1.08309592629463; "I" = 1.049984713448011; "J" = 1.075718823968999; "K" = 1.074702621514289; "L" = 1.076747725040804; "M" = 1.086631764718973; "N" = 1.077246468085709 }
    23 = @{ "B" = 1.02; "C" = 1.069820814771376; "D" = 1.071680565473826; "E" = 1.073795284062127; "F" = 1.083729795163096; "I" = 1.050186878054584; "J" = 1.076282284256428; "K" = 1.075191172887073; "L" = 1.077298472864302; "M" = 1.087198559237673; "N" = 1.077810728551376 }
    24 = @{ "B" = 1.02; "C" = 1.072538535722966; "D" = 1.073863184100487; "E" = 1.076222449778967; "F" = 1.08621812332631; "I" = 1.050975922236093; "J" = 1.078492054961257; "K" = 1.077106224555616; "L" = 1.079457899584747; "M" = 1.089421813285632; "N" = 1.08002363738407 }
    25 = @{ "B" = 1.02; "C" = 1.075676384960971; "D" = 1.076381562164344; "E" = 1.079023509569926; "F" = 1.089091487906264; "I" = 1.051877595156984; "J" = 1.081039313273928; "K" = 1.079311788638867; "L" = 1.081946120671694; "M" = 1.091985445823412; "N" = 1.082574513095719 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
